# New file with texts
#
# - Adds a new shared string "Follow" as a new row on the "Basic" sheet,
#   right before the existing "Translate the extension" row (old row 48),
#   pushing every row below it down by one.
# - Makes "More" the active/selected sheet (was "Basic").
# - Updates each sheet's scrolled position to match where the user had
#   scrolled to when the workbook was saved.

$wb = $excel.ActiveWorkbook
$basic = $wb.Worksheets.Item("Basic")
$more = $wb.Worksheets.Item("More")

# Insert a new row above the old row 48 and give it the new text "Follow".
# This shifts rows 48:81 down to 49:82, carrying their styles/content with them.
$basic.Rows.Item(48).Insert()
$basic.Cells.Item(48, 1).Value = "Follow"

# Restore the selected cell on "Basic" (unchanged by the diff) and scroll it
# so row 67 is at the top of the window.
$basic.Activate()
$basic.Range("A49").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1

# "More" becomes the active sheet/tab, scrolled back to the top, with its
# own previous selection kept.
$more.Activate()
$more.Range("A25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
